$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95 - Primera
$ws.Cells.Item(95, 1).Value = 1
$ws.Cells.Item(95, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(95, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(95, 4).Value = 44628
$ws.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(95, 5).Value = 15
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100108
$ws.Cells.Item(95, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(95, 9).Value = 100108003
$ws.Cells.Item(95, 10).Value = "Maracuyá"
$ws.Cells.Item(95, 11).Value = "Sin especificar"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 120
$ws.Cells.Item(95, 14).Value = 18000
$ws.Cells.Item(95, 15).Value = 19000
$ws.Cells.Item(95, 16).Value = 18500
$ws.Cells.Item(95, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(95, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(95, 19).Value = 925
$ws.Cells.Item(95, 20).Value = 20

# Row 96 - Segunda
$ws.Cells.Item(96, 1).Value = 1
$ws.Cells.Item(96, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(96, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(96, 4).Value = 44628
$ws.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(96, 5).Value = 15
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100108
$ws.Cells.Item(96, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(96, 9).Value = 100108003
$ws.Cells.Item(96, 10).Value = "Maracuyá"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Segunda"
$ws.Cells.Item(96, 13).Value = 100
$ws.Cells.Item(96, 14).Value = 17000
$ws.Cells.Item(96, 15).Value = 18000
$ws.Cells.Item(96, 16).Value = 17500
$ws.Cells.Item(96, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(96, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96, 19).Value = 875
$ws.Cells.Item(96, 20).Value = 20
